# Sincronização de dados: adiciona novo orçamento (quotations row 18) e os
# respectivos itens de orçamento (items rows 55-66).

$wb = $excel.ActiveWorkbook
$wsQuotations = $wb.Worksheets.Item("quotations")
$wsItems = $wb.Worksheets.Item("items")

# Helper: write a value that must stay numeric-like text (e.g. "549.22", "0")
# as a real text cell instead of letting Excel auto-convert it to a number.
function Set-TextValue($sheet, $row, $col, $text) {
    if ($text -eq $null) { return }
    $sheet.Cells.Item($row, $col).Value = "'" + $text
}

# Helper: write a plain value (string that is not numeric-looking, or a
# boolean/number) straight through.
function Set-RawValue($sheet, $row, $col, $value) {
    if ($value -eq $null) { return }
    $sheet.Cells.Item($row, $col).Value = $value
}

# ---------------------------------------------------------------------------
# quotations!A18:U18 — new quotation row
# ---------------------------------------------------------------------------

Set-RawValue  $wsQuotations 18 1  "ZWQ5MzRmMGYtZDE1Yy00ODY1LWFjMDAtZTY4ZWE1YzRkODk5OjU3MDE2"
Set-RawValue  $wsQuotations 18 2  "R1BQOWO2EE"
Set-RawValue  $wsQuotations 18 3  "GRALHA AZUL PONTA GROSSA"
Set-RawValue  $wsQuotations 18 4  "ORÇAMENTO REFERENTE BALDE DANIFICADO POR MAL USO`nBALDE NOVO , SE FOR ENTREGAR O BALDE 43,00 PELO CORREIO VIA SEDEX"
# E18 -> additionalInformation: left blank
Set-RawValue  $wsQuotations 18 6  $false
Set-TextValue $wsQuotations 18 7  "549.22"
Set-TextValue $wsQuotations 18 8  "549.22"
Set-RawValue  $wsQuotations 18 9  "Pendente"
Set-RawValue  $wsQuotations 18 10 "2025-10-06T19:05:55.326Z"
# K18 -> refusedAt: left blank
# L18 -> approvedBy: left blank
Set-RawValue  $wsQuotations 18 13 "Adriana Vieira Masini"
# N18 -> approvedSignature: left blank
Set-RawValue  $wsQuotations 18 15 "2025-09-29T20:26:53.403Z"
Set-RawValue  $wsQuotations 18 16 "MDBmZTJkNzktN2M1YS00MDc0LWE2YTctNzZiZGNkZWFmYTIwOjU3MDE2"
Set-RawValue  $wsQuotations 18 17 "percentage"
Set-TextValue $wsQuotations 18 18 "0"
Set-TextValue $wsQuotations 18 19 "0"
Set-RawValue  $wsQuotations 18 20 "NDgyNTE1Mzo1NzAxNg=="
Set-RawValue  $wsQuotations 18 21 "pending"

# The multi-line description above makes Excel auto-grow the row height;
# reset it so row 18 keeps the sheet's default height like every other row.
$wsQuotations.Rows.Item(18).AutoFit()

# ---------------------------------------------------------------------------
# items!A55:J66 — new item rows tied to the quotation above
# ---------------------------------------------------------------------------

$quotationRef = "ZWQ5MzRmMGYtZDE1Yy00ODY1LWFjMDAtZTY4ZWE1YzRkODk5OjU3MDE2"

$items = @(
    @{ Row=55; A="MWUxOTRiYWItMTljOC00NTMwLWIzODQtZjA3YTRjMDYzMjZlOjU3MDE2"; C=43;              D="CORREIO VIA SEDEX 43,00"; G="MTZhMTdiZDUtYWI4OC00NjYwLTk2OGUtOWE4YmRiY2JmYTcyOjU3MDE2"; H=43 },
    @{ Row=56; A="MzViYTlhNTItOTE2MC00NmJlLTlkNGMtN2MyZTQ3M2NlNzFiOjU3MDE2"; C=3325;            D=$null;                     G="ZDNmZWYzYTctMmQwZC00YTJjLWJjYjAtYTZiYzQ4ZjBiYTIxOjU3MDE2"; H=3325 },
    @{ Row=57; A="NGNmZDE3ZWQtZTFhMy00MWZkLWE3ZTgtYjg5ZDA1MTQ2MjUwOjU3MDE2"; C=717;             D=$null;                     G="OWQzZDUyMzctNGVhYS00ZDcwLWIxYzQtNjJjM2VjYTEzYjJmOjU3MDE2"; H=717 },
    @{ Row=58; A="NmRhZDJlMTQtZWQ0My00NWIzLThhMmMtYmM1OTIxYTUxN2UwOjU3MDE2"; C=7375;            D=$null;                     G="ZGRlNDk5MWMtYjg3Ny00N2MzLWE2MDgtNjc3MWVhODJjYjlhOjU3MDE2"; H=7375 },
    @{ Row=59; A="OGZmOWRhOTQtNDdiZC00MDNlLWE5MTctZTc5MTVlY2FmZWYzOjU3MDE2"; C=117;             D=$null;                     G="YTQ3NzcxODMtYTBkNS00ZWUzLWJlYTQtYjMyNzI5MTRmODhhOjU3MDE2"; H=117 },
    @{ Row=60; A="YWQ4YTM5MmEtZDQ0OC00ZjdhLWE3NjQtNGI0MjJlYzBjNGFkOjU3MDE2"; C=42000000000000000; D=$null;                     G="NGRmZWRlNWQtYjFlNi00N2Q2LTkxZDItMDBmZjM2MzVhMTUxOjU3MDE2"; H=42000000000000000 },
    @{ Row=61; A="YWVmZTk0M2ItZWU2ZS00ODMyLTgzZWUtMGYwMTdkMmVjMmYwOjU3MDE2"; C=6449999999999999;  D=$null;                     G="MGFmNjg5ODYtODc1ZC00YmFjLWE1MDgtMTE3YWU5MjExMzVkOjU3MDE2"; H=6449999999999999 },
    @{ Row=62; A="YjA5NTJjYmItOTU2MS00NThlLTliZDQtMzhkNGQwNWY4MGFmOjU3MDE2"; C=48;              D=$null;                     G="ODcwZTI1ZDEtMTRkNC00M2IyLTk0MTItOGJhNDdiYzIzMjg1OjU3MDE2"; H=48 },
    @{ Row=63; A="Yjk1N2RlZTctZDMxZi00YzJlLWI1NTktZmEyMjQ1NjQ1ZDg3OjU3MDE2"; C=475;             D=$null;                     G="ODNiOGMzNWItYjY5NC00ZWNlLWFjYzAtNGFmYjE5MTRjZjY0OjU3MDE2"; H=475 },
    @{ Row=64; A="YmY5ZDcxOGUtOTg3Ny00YjYyLTk4ZmMtZDY4MTk3OTM4MmFjOjU3MDE2"; C=1725;            D=$null;                     G="ODJmMjkxYjEtMDdiOC00YjFiLWEyZmMtMjExYjg2YjIwYjgyOjU3MDE2"; H=1725 },
    @{ Row=65; A="YzMxNzM5OTQtY2Y0MS00NjdmLWI3MzctMzVhYTk3ODgzMWM4OjU3MDE2"; C=6000000000000001;  D=$null;                     G="OTgxYjZlMTAtNGZiMy00YjAwLWI4OTYtMTcxNGM5MTg2Y2NiOjU3MDE2"; H=6000000000000001 },
    @{ Row=66; A="ZTNmNTk0NDctYjYxYy00MzgyLWE0NWQtZDY5NjI3NTdkNGQ4OjU3MDE2"; C=8399999999999999;  D=$null;                     G="MmY2Nzk2OGUtMDNlOC00MzZlLTk0ZjAtZTE1ODNmMDk3OTJiOjU3MDE2"; H=8399999999999999 }
)

foreach ($item in $items) {
    $r = $item.Row
    Set-RawValue $wsItems $r 1  $item.A          # item_id
    Set-RawValue $wsItems $r 2  1                # item_quantity
    Set-RawValue $wsItems $r 3  $item.C           # item_total
    Set-RawValue $wsItems $r 4  $item.D           # item_description
    Set-RawValue $wsItems $r 5  13                # item_position
    Set-RawValue $wsItems $r 6  $quotationRef      # item_quotation_ref
    Set-RawValue $wsItems $r 7  $item.G           # productService_id
    Set-RawValue $wsItems $r 8  $item.H           # productService_value
    Set-RawValue $wsItems $r 9  "product"          # productService_type
    Set-RawValue $wsItems $r 10 $quotationRef      # quotation_id
}
